# Reproduce the commit: "Added new methods to ProfilePage and PostRecordViewPage"
# New shared strings (must be introduced in this exact order so they land at
# sharedStrings indices 152/153/154, matching the target OOXML):
#   152 CreateAndEditPost
#   153 OPQA-382|OPQA-388|OPQA-406|OPQA-372
#   154 Verify that user is able to create and edit post and verify that time stamp is displayed

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Test Cases" (sheet1) - new row 41 with a new test case record.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Test Cases")
$ws1.Activate()

$ws1.Range("A40:E40").Copy()
$ws1.Range("A41:E41").PasteSpecial(-4122)

$ws1.Range("A41").Value = "CreateAndEditPost"
$ws1.Range("B41").Value = "OPQA-382|OPQA-388|OPQA-406|OPQA-372"
$ws1.Range("C41").Value = "Verify that user is able to create and edit post and verify that time stamp is displayed"
$ws1.Range("D41").Value = "Y"
$ws1.Range("E41").Value = "PASS"

$ws1.Hyperlinks.Add($ws1.Range("B41"), "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-382", "", "", "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-382")

$ws1.Range("A43:B44").Select()

# ---------------------------------------------------------------------------
# Sheet "Test Case Steps" (sheet2) - extend used range to column E / row 41
# and record the same step values.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Test Case Steps")
$ws2.Activate()

$ws2.Range("A2").Copy()
$ws2.Range("D1:E1").PasteSpecial(-4122)
$ws2.Range("D2:E4").PasteSpecial(-4122)

$ws2.Range("A2:E2").Copy()
$ws2.Range("A5:E40").PasteSpecial(-4122)
$ws2.Range("A41:E41").PasteSpecial(-4122)

$ws2.Range("B41").Value = "OPQA-382|OPQA-388|OPQA-406|OPQA-372"
$ws2.Range("C41").Value = "Verify that user is able to create and edit post and verify that time stamp is displayed"
$ws2.Range("D41").Value = "Y"

$ws2.Range("A43:B44").Select()

# ---------------------------------------------------------------------------
# Sheet "AuthoringTest" (sheet3)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("AuthoringTest")
$ws3.Activate()

$ws3.Range("A2:E2").Copy()
$ws3.Range("A3:E40").PasteSpecial(-4122)
$ws3.Range("A41:E41").PasteSpecial(-4122)

$ws3.Range("B41").Value = "OPQA-382|OPQA-388|OPQA-406|OPQA-372"
$ws3.Range("C41").Value = "Verify that user is able to create and edit post and verify that time stamp is displayed"
$ws3.Range("D41").Value = "Y"

$ws3.Range("A43:B44").Select()

# ---------------------------------------------------------------------------
# Sheet "CommentsMinMaxValidationTest" (sheet4)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("CommentsMinMaxValidationTest")
$ws4.Activate()

$ws4.Range("A2:E2").Copy()
$ws4.Range("A3:E40").PasteSpecial(-4122)
$ws4.Range("A41:E41").PasteSpecial(-4122)

$ws4.Range("B41").Value = "OPQA-382|OPQA-388|OPQA-406|OPQA-372"
$ws4.Range("C41").Value = "Verify that user is able to create and edit post and verify that time stamp is displayed"
$ws4.Range("D41").Value = "Y"

$ws4.Range("A43:B44").Select()

# ---------------------------------------------------------------------------
# Sheet "CommentsProfanityWordsCheckTest" (sheet5)
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("CommentsProfanityWordsCheckTest")
$ws5.Activate()

$ws5.Range("A2:D2").Copy()
$ws5.Range("E1:E7").PasteSpecial(-4122)
$ws5.Range("A8:E40").PasteSpecial(-4122)
$ws5.Range("A41:E41").PasteSpecial(-4122)

$ws5.Range("B41").Value = "OPQA-382|OPQA-388|OPQA-406|OPQA-372"
$ws5.Range("C41").Value = "Verify that user is able to create and edit post and verify that time stamp is displayed"
$ws5.Range("D41").Value = "Y"

$ws5.Range("A43:B44").Select()

$ws1.Activate()
